# Update metadata: reclassify "horas-trabajadas" and "sector-actividad-descripcion"
# from measures to dimensions, and add mapping file references.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: semantic identifiers - measure -> dimension for columns A and C
$ws.Range("A3").Value = "iaest-dimension:horas-trabajadas"
$ws.Range("C3").Value = "iaest-dimension:sector-actividad-descripcion"

# Row 4: category - medida -> dim for columns A and C
$ws.Range("A4").Value = "dim"
$ws.Range("C4").Value = "dim"

# Row 5: datatype - xsd:string -> skos:Concept for columns A and C
$ws.Range("A5").Value = "skos:Concept"
$ws.Range("C5").Value = "skos:Concept"

# Row 6 (new): mapping file references for the newly-reclassified dimensions.
# Copy formatting from row 5 so the new cells share the same style as the rest
# of the table, then set their values.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = "mapping-horas-trabajadas.xlsx"
$ws.Range("C6").Value = "mapping-sector-actividad-descripcion.xlsx"
